# Add a new OJT log entry (row 2) to the tracker sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format A2 as Text first so the date-like string "2025-07-29" is kept
# as literal text instead of being auto-converted into a date serial
# number by Excel's input parser.
$ws.Range("A2").NumberFormat = "@"

$ws.Range("A2").Value = "2025-07-29"
$ws.Range("B2").Value = "Tuesday"
$ws.Range("C2").Value = "08:30"
$ws.Range("D2").Value = "17:30"
$ws.Range("E2").Value = 9
